$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7522
$ws1.Range("F5").Value = 13
$ws1.Range("F7").Value = 4121
$ws1.Range("F9").Value = 577
$ws1.Range("F11").Value = 660
$ws1.Range("F12").Value = 147

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7522
$ws4.Range("F7").Value = 13
$ws4.Range("F9").Value = 4121
$ws4.Range("F11").Value = 577
$ws4.Range("F13").Value = 660
$ws4.Range("F15").Value = 147
